$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 74
$ws.Range("H74").Value = 4000
$ws.Range("I74").Value = 0
$ws.Range("K74").Value = 0
$ws.Range("M74").Value = $null

# Row 77
$ws.Range("H77").Value = 4000
$ws.Range("I77").Value = 0
$ws.Range("K77").Value = 0
$ws.Range("M77").Value = $null

# Row 88
$ws.Range("H88").Value = 1235558.6
$ws.Range("J88").Value = 1764920.9
$ws.Range("L88").Value = 1764920.9
$ws.Range("N88").Value = -1765732.9

# Row 91
$ws.Range("H91").Value = 1235558.6
$ws.Range("J91").Value = 1764920.9
$ws.Range("L91").Value = 1764920.9
$ws.Range("N91").Value = -1767728.9

# Row 112
$ws.Range("H112").Value = 2128.6099
$ws.Range("J112").Value = 2408.647
$ws.Range("L112").Value = 7225.941
$ws.Range("N112").Value = -9441.940999999999

# Row 113
$ws.Range("H113").Value = 3001.8333
$ws.Range("I113").Value = 2502.5
$ws.Range("J113").Value = 3251.5
$ws.Range("K113").Value = 2502.5
$ws.Range("L113").Value = 3251.5
$ws.Range("M113").Value = 751.5
$ws.Range("N113").Value = -9759.5

# Row 116
$ws.Range("H116").Value = 3600.125
$ws.Range("I116").Value = 3300.0667
$ws.Range("K116").Value = 3300.0667
$ws.Range("M116").Value = 141.9333000000001

# Row 129
$ws.Range("H129").Value = 846.3099999999999
$ws.Range("I129").Value = 329.5
$ws.Range("J129").Value = 891.25
$ws.Range("K129").Value = 988.5
$ws.Range("L129").Value = 2673.75
$ws.Range("M129").Value = 4011.5
$ws.Range("N129").Value = -12673.75

# Row 134
$ws.Range("H134").Value = 0
$ws.Range("J134").Value = 0
$ws.Range("L134").Value = 0
$ws.Range("N134").Value = $null

# Row 137
$ws.Range("H137").Value = 1720.4375
$ws.Range("I137").Value = 1021.2
$ws.Range("J137").Value = 2038.2727
$ws.Range("K137").Value = 3063.6
$ws.Range("L137").Value = 6114.8181
$ws.Range("M137").Value = -513.6000000000004
$ws.Range("N137").Value = -11214.8181

# Row 138
$ws.Range("H138").Value = 1451.8889
$ws.Range("I138").Value = 614.7805
$ws.Range("J138").Value = 2043.638
$ws.Range("K138").Value = 1844.3415
$ws.Range("L138").Value = 6130.914
$ws.Range("M138").Value = 3295.6585
$ws.Range("N138").Value = -16410.914

$ws = $wb.Worksheets.Item("ARM")
# Row 32
$ws.Range("H32").Value = 2580.3289
$ws.Range("I32").Value = 2578.9077
$ws.Range("J32").Value = 2588.7273
$ws.Range("K32").Value = 2578.9077
$ws.Range("L32").Value = 2588.7273
$ws.Range("M32").Value = -2291.9077
$ws.Range("N32").Value = -3162.7273

# Row 92
$ws.Range("H92").Value = 5000000
$ws.Range("J92").Value = 5000000
$ws.Range("L92").Value = 5000000
$ws.Range("N92").Value = -5004992

# Row 132
$ws.Range("H132").Value = 1145.6938
$ws.Range("I132").Value = 866.06976
$ws.Range("K132").Value = 2598.20928
$ws.Range("M132").Value = -68.20928000000004

$ws = $wb.Worksheets.Item("BSM")
# Row 107
$ws.Range("H107").Value = 1636.3529
$ws.Range("I107").Value = 1481.2
$ws.Range("J107").Value = 1858
$ws.Range("K107").Value = 1481.2
$ws.Range("L107").Value = 1858
$ws.Range("M107").Value = 438.8
$ws.Range("N107").Value = -5698

# Row 134
$ws.Range("H134").Value = 3075.5356
$ws.Range("I134").Value = 878.6667
$ws.Range("J134").Value = 12062.728
$ws.Range("K134").Value = 2636.0001
$ws.Range("L134").Value = 36188.18399999999
$ws.Range("M134").Value = -101.0001000000002
$ws.Range("N134").Value = -41258.18399999999

$ws = $wb.Worksheets.Item("CRP")
# Row 31
$ws.Range("H31").Value = 1198.2407
$ws.Range("I31").Value = 1174.22
$ws.Range("J31").Value = 1498.5
$ws.Range("K31").Value = 1174.22
$ws.Range("L31").Value = 1498.5
$ws.Range("M31").Value = -879.22
$ws.Range("N31").Value = -2088.5

# Row 34
$ws.Range("H34").Value = 1198.2407
$ws.Range("I34").Value = 1174.22
$ws.Range("J34").Value = 1498.5
$ws.Range("K34").Value = 1174.22
$ws.Range("L34").Value = 1498.5
$ws.Range("M34").Value = -972.22
$ws.Range("N34").Value = -1902.5

# Row 43
$ws.Range("H43").Value = 4445
$ws.Range("J43").Value = 4445
$ws.Range("L43").Value = 4445
$ws.Range("N43").Value = -4813

# Row 99
$ws.Range("H99").Value = 2633334
$ws.Range("I99").Value = 3761048.8
$ws.Range("K99").Value = 3761048.8
$ws.Range("M99").Value = -3759550.8

# Row 101
$ws.Range("H101").Value = 4445
$ws.Range("J101").Value = 4445
$ws.Range("L101").Value = 4445
$ws.Range("N101").Value = -10935

# Row 106
$ws.Range("H106").Value = 49899.5
$ws.Range("J106").Value = 49899.5
$ws.Range("L106").Value = 49899.5
$ws.Range("N106").Value = -52423.5

# Row 126
$ws.Range("H126").Value = 2633334
$ws.Range("I126").Value = 3761048.8
$ws.Range("K126").Value = 11283146.4
$ws.Range("M126").Value = -11280676.4

$ws = $wb.Worksheets.Item("CUL")
# Row 14
$ws.Range("H14").Value = 354
$ws.Range("I14").Value = 354
$ws.Range("K14").Value = 1062
$ws.Range("M14").Value = -889

# Row 131
$ws.Range("H131").Value = 13335332
$ws.Range("J131").Value = 2188.597
$ws.Range("L131").Value = 6565.791000000001
$ws.Range("N131").Value = -16645.791

# Row 136
$ws.Range("H136").Value = 1347.5
$ws.Range("I136").Value = 922.3077
$ws.Range("J136").Value = 2453
$ws.Range("K136").Value = 2766.9231
$ws.Range("L136").Value = 7359
$ws.Range("M136").Value = 2333.0769
$ws.Range("N136").Value = -17559

$ws = $wb.Worksheets.Item("GSM")
# Row 126
$ws.Range("H126").Value = 2314.0715
$ws.Range("I126").Value = 1799.625
$ws.Range("K126").Value = 5398.875
$ws.Range("M126").Value = -2928.875

# Row 127
$ws.Range("H127").Value = 34352.94
$ws.Range("J127").Value = 34352.94
$ws.Range("L127").Value = 34352.94
$ws.Range("N127").Value = -44272.94

# Row 132
$ws.Range("H132").Value = 1886.2554
$ws.Range("I132").Value = 1348.0714
$ws.Range("K132").Value = 4044.2142
$ws.Range("M132").Value = -1514.2142

$ws = $wb.Worksheets.Item("LTW")
# Row 93
$ws.Range("H93").Value = 758.5714
$ws.Range("I93").Value = 758.5714
$ws.Range("J93").Value = 0
$ws.Range("K93").Value = 758.5714
$ws.Range("L93").Value = 0
$ws.Range("M93").Value = 489.4286
$ws.Range("N93").Value = $null

# Row 100
$ws.Range("H100").Value = 1369.6666
$ws.Range("I100").Value = 804.5
$ws.Range("K100").Value = 804.5
$ws.Range("M100").Value = -263.5

# Row 101
$ws.Range("H101").Value = 17000
$ws.Range("J101").Value = 17000
$ws.Range("L101").Value = 17000
$ws.Range("N101").Value = -23490

# Row 106
$ws.Range("H106").Value = 16250
$ws.Range("J106").Value = 16250
$ws.Range("L106").Value = 16250
$ws.Range("N106").Value = -18774

# Row 132
$ws.Range("H132").Value = 27981.525
$ws.Range("I132").Value = 1045.2273
$ws.Range("K132").Value = 3135.6819
$ws.Range("M132").Value = -605.6819

# Row 134
$ws.Range("H134").Value = 28615.455
$ws.Range("J134").Value = 28615.455
$ws.Range("L134").Value = 28615.455
$ws.Range("N134").Value = -38755.455

# Row 136
$ws.Range("H136").Value = 1855.8182
$ws.Range("I136").Value = 2251
$ws.Range("J136").Value = 1630
$ws.Range("K136").Value = 6753
$ws.Range("L136").Value = 4890
$ws.Range("M136").Value = -4203
$ws.Range("N136").Value = -9990

$ws = $wb.Worksheets.Item("WVR")
# Row 92
$ws.Range("H92").Value = 14587.5
$ws.Range("J92").Value = 14587.5
$ws.Range("L92").Value = 14587.5
$ws.Range("N92").Value = -19579.5

# Row 132
$ws.Range("H132").Value = 1117.4348
$ws.Range("I132").Value = 679.1053000000001
$ws.Range("J132").Value = 3199.5
$ws.Range("K132").Value = 2037.3159
$ws.Range("L132").Value = 9598.5
$ws.Range("M132").Value = 492.6840999999999
$ws.Range("N132").Value = -14658.5

# Row 136
$ws.Range("H136").Value = 1032.125
$ws.Range("I136").Value = 653
$ws.Range("J136").Value = 1411.25
$ws.Range("K136").Value = 1959
$ws.Range("L136").Value = 4233.75
$ws.Range("M136").Value = 591
$ws.Range("N136").Value = -9333.75
